$p = $ppt.ActivePresentation

# The deck grows from 8 slides to 11:
#   1 Title                 (unchanged)
#   2 Short Overview        (new)
#   3 Introduction          (unchanged)
#   4 Requirements          (unchanged)
#   5 Class Diagram         (unchanged)
#   6 Use-Case Diagram      (unchanged)
#   7 Activity Diagram      (unchanged)
#   8 Sequence Diagram      (unchanged)
#   9 Prototypes            (new)
#  10 Review and Conclusion (was "Conclusion", text edited in place)
#  11 Thank You             (new)
#
# New slides are created (duplicating the "Introduction" slide as a
# structural template so placeholders/formatting match the rest of the deck)
# in the same chronological order PowerPoint allocated their internal slide
# ids in the source edit: Prototypes first, then Short Overview, then
# Thank You last - even though their final on-slide positions differ.

# 1) Template for "Prototypes" - duplicate "Introduction" (index 2) and move
#    it to sit right before "Conclusion" (currently the last slide, index 8).
$dupPrototypes = $p.Slides.Item(2).Duplicate()
$dupPrototypes.Item(1).MoveTo(8)

# 2) Template for "Short Overview" - duplicate "Introduction" (still index 2)
#    and move it right after the title slide.
$dupShortOverview = $p.Slides.Item(2).Duplicate()
$dupShortOverview.Item(1).MoveTo(2)

# 3) Fill in the titles for the two new slides.
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Text = "Short Overview"
$p.Slides.Item(9).Shapes.Item(1).TextFrame.TextRange.Text = "Prototypes"

# 4) Rename "Conclusion" (now at index 10) to "Review and Conclusion" - same
#    slide identity, only the title text changes.
$p.Slides.Item(10).Shapes.Item(1).TextFrame.TextRange.Text = "Review and Conclusion"

# 5) Template for "Thank You" - duplicate "Introduction" (still index 3) and
#    append it at the very end of the deck.
$dupThankYou = $p.Slides.Item(3).Duplicate()
$dupThankYou.Item(1).MoveTo($p.Slides.Count)
$p.Slides.Item($p.Slides.Count).Shapes.Item(1).TextFrame.TextRange.Text = "Thank You"
